$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Coefficients")
$ws1.Range("A5").Value = "MonthOctober"
$ws1.Range("A7").Value = "DRM:MonthOctober"
$ws1.Range("A8").Value = "Depth:MonthOctober"

$ws2 = $wb.Worksheets.Item("Fullmodel_statistics")
$ws2.Range("A2").Value = 0.6518858270538
$ws2.Range("B2").Value = 0.5715517871431385
$ws2.Range("C2").Value = 8.114689959309327
$ws2.Range("E2").Value = 0.0000530807759883436
